$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.13"
$ws.Range("E2").Value = "'1.23%"
$ws.Range("D3").Value = "'41.09"
$ws.Range("E3").Value = "'1.87%"
$ws.Range("D4").Value = "'5.117"
$ws.Range("E4").Value = "'0.27%"
$ws.Range("D5").Value = "'0.07690"
$ws.Range("E5").Value = "'1.40%"
$ws.Range("D6").Value = "'4.286"
$ws.Range("E6").Value = "'0.24%"
$ws.Range("D7").Value = "'1.619"
$ws.Range("E7").Value = "'0.68%"
$ws.Range("D8").Value = "'0.9202"
$ws.Range("E8").Value = "'1.72%"
$ws.Range("D10").Value = "'0.1220"
$ws.Range("E10").Value = "'20.83%"
$ws.Range("D11").Value = "'0.1836"
$ws.Range("E11").Value = "'4.70%"
$ws.Range("D12").Value = "'0.09101"
$ws.Range("E12").Value = "'0.14%"
$ws.Range("D13").Value = "'0.04257"
$ws.Range("E13").Value = "'0.82%"
$ws.Range("E14").Value = "'-0.33%"
$ws.Range("D15").Value = "'0.001244"
$ws.Range("E15").Value = "'1.47%"
$ws.Range("D16").Value = "'0.005828"
$ws.Range("E16").Value = "'-0.04%"
$ws.Range("D17").Value = "'3.353"
$ws.Range("E17").Value = "'0.14%"
$ws.Range("E18").Value = "'1.22%"
$ws.Range("D19").Value = "'6.906"
$ws.Range("E19").Value = "'3.74%"
$ws.Range("D20").Value = "'0.1379"
$ws.Range("E20").Value = "'1.48%"
$ws.Range("D21").Value = "'0.2675"
$ws.Range("E21").Value = "'-2.06%"
$ws.Range("D22").Value = "'0.04035"
$ws.Range("E22").Value = "'-3.58%"
$ws.Range("D23").Value = "'0.001262"
$ws.Range("E23").Value = "'2.63%"
$ws.Range("D24").Value = "'0.004076"
$ws.Range("E24").Value = "'0.67%"
$ws.Range("D25").Value = "'0.0001267"
$ws.Range("E25").Value = "'-2.72%"
$ws.Range("E26").Value = "'24.48%"
$ws.Range("D38").Value = "'0.02474"
$ws.Range("E38").Value = "'3.69%"
$ws.Range("D39").Value = "'0.05267"
$ws.Range("E39").Value = "'2.62%"
$ws.Range("D40").Value = "'0.007828"
$ws.Range("E40").Value = "'0.68%"
$ws.Range("D41").Value = "'0.1313"
$ws.Range("E41").Value = "'1.47%"
$ws.Range("D42").Value = "'0.006790"
$ws.Range("E42").Value = "'-3.76%"
$ws.Range("E43").Value = "'-5.40%"
$ws.Range("D44").Value = "'0.008184"
$ws.Range("E44").Value = "'-3.32%"
$ws.Range("E45").Value = "'-6.61%"
$ws.Range("D46").Value = "'0.00006838"
$ws.Range("E46").Value = "'7.23%"
$ws.Range("D47").Value = "'0.00000000748"
$ws.Range("E47").Value = "'-0.42%"
$ws.Range("D48").Value = "'0.2286"
$ws.Range("E48").Value = "'3,082.11%"
$ws.Range("D49").Value = "'0.004091"
$ws.Range("E49").Value = "'-7.20%"
$ws.Range("E50").Value = "'-0.42%"
$ws.Range("E51").Value = "'-0.42%"
